# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 -------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 450
$wsExhibit.Range("F3").Value  = 32
$wsExhibit.Range("F4").Value  = 62
$wsExhibit.Range("F5").Value  = 5080
$wsExhibit.Range("F6").Value  = 172
$wsExhibit.Range("F7").Value  = 21
$wsExhibit.Range("F9").Value  = 309
$wsExhibit.Range("F10").Value = 55

# --- Sheet: 全部类型 -----------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 450
$wsAll.Range("F7").Value  = 32
$wsAll.Range("F8").Value  = 62
$wsAll.Range("F9").Value  = 5080
$wsAll.Range("F10").Value = 172
$wsAll.Range("F11").Value = 22
$wsAll.Range("F14").Value = 309
$wsAll.Range("F15").Value = 55
